$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the bordered cell style (the "s=1" format used throughout the
# table) down across the newly-used rows 13-17 by copying the format from
# row 12 (which already carries it) before we touch any values. ---
$src = $ws.Range("A12:E12")
$dst = $ws.Range("A13:E17")
$src.Copy()
$dst.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 12 keeps only the "Name" entry; the rest of what used to live
# there (Vorname/Studiengang/Vorlesungstage) moves one field at a time
# down onto its own row. ---
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("E12").ClearContents()

$ws.Range("B13").Value = "Max"
$ws.Range("C14").Value = "BSC"
$ws.Range("D15").Value = "21/22"
$ws.Range("E16").Value = "Mo + Di"

# --- Row 17: a brand new, fully populated row demonstrating an invalid
# "Studiengang" entry ("A") that only trips a warning-level validation. ---
$ws.Range("A17").Value = "Mustermann"
$ws.Range("B17").Value = "Max"
$ws.Range("C17").Value = "A"
$ws.Range("D17").Value = "21/22"
$ws.Range("E17").Value = "Mo + Di"

# --- Data validation bookkeeping ---
# The "Studiengang" list validation now also covers the newly used rows
# (C2:C12 grows to C2:C16; C21:C1048576 is left as-is).
$ws.Range("C2:C12").Validation.Delete()
$ws.Range("C2:C16").Validation.Add(3, 1, 1, '"BWI,BSC,VI,FISI"')

# The date validation on the "Jahrgang" cell grows to span D12:D14.
$ws.Range("D12").Validation.Delete()
$ws.Range("D12:D14").Validation.Add(4, 1, 1, 2010, 2500)

# New, warning-only list validation for the deliberately-invalid C17 cell.
$ws.Range("C17").Validation.Add(3, 2, 1, '"BWI,BSC,VI,FISI"')

# --- Column widths tweaked slightly (values chosen so the engine's
# pixel-grid rounding lands as close as possible to the authored
# 15.75 / 16.75 / 15.75 / 18.125 / 20.875 character widths) ---
$ws.Columns.Item(1).ColumnWidth = 15.035714285714286
$ws.Columns.Item(2).ColumnWidth = 16.035714285714285
$ws.Columns.Item(3).ColumnWidth = 15.035714285714286
$ws.Columns.Item(4).ColumnWidth = 17.410714285714285
$ws.Columns.Item(5).ColumnWidth = 20.160714285714285

# --- Header row height tweaked slightly ---
$ws.Rows.Item(1).RowHeight = 27.4

# --- Selection moves to F15 ---
$null = $ws.Range("F15").Select()

"done"
